# Weekly update: insert a new latest-week record for
# "Hortaliza, Agrícola del Norte S.A. de Arica - Cebollín baby" at row 23,
# pushing the previously-existing rows 23:52 down to 24:53.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 23 (shifts rows 23:52 -> 24:53, dimension grows to R53).
$ws.Rows("23:23").Insert()

# Populate the new row with this week's record (same market / product metadata
# as every other row in this sheet, new date + price figures).
$ws.Range("A23").Value = 1
$ws.Range("B23").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C23").Value = "Arica y Parinacota"
$ws.Range("D23").Value2 = 44495
$ws.Range("E23").Value = 15
$ws.Range("F23").Value = 100112038
$ws.Range("G23").Value = "Cebollín baby"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 1000
$ws.Range("L23").Value = 1200
$ws.Range("M23").Value = 1100
$ws.Range("N23").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O23").Value = "Región de Arica y Parinacota"
$ws.Range("P23").Value = 550
$ws.Range("Q23").Value = 2
$ws.Range("R23").Value = "Hortaliza"
